# Mise en fonctionnement de la page interfaceCaracteristique.php
#
# 1. Rename the shared group label "Loire" -> "Plaine" everywhere it is used
#    (column A on both sheets).
# 2. Update the "Nombre d'exploitations" (D), "Nombre de mesures" (E) and
#    "Moyennes croissance" (F) figures for a batch of weeks/décades on both
#    the "2021 Semaine" and "2021 Décade" sheets - some rows gain figures
#    that were previously blank, some rows' figures are cleared out, and
#    most existing rows get revised numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: rename "Loire" to "Plaine" on every worksheet.
# ---------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Loire", "Plaine") | Out-Null
}

# ---------------------------------------------------------------------
# Step 2: per-row D/E/F updates.
# Each tuple is (row, D-value, E-value, F-value).
# "SKIP"  -> leave the existing cell untouched
# "CLEAR" -> blank the cell out
# ---------------------------------------------------------------------

$sheet1Updates = @(
    @(9,  'CLEAR', 'CLEAR', 'CLEAR'),
    @(11, '2',  '21',  '14.15'),
    @(12, '4',  '55',  '20.18'),
    @(13, 'SKIP', '56', '24.31'),
    @(14, '10', '166', '32.94'),
    @(15, '7',  '74',  '44.95'),
    @(16, '12', '112', '42.3'),
    @(17, '14', '127', '34.39'),
    @(18, '11', '108', '49.96'),
    @(19, '14', '110', '69.17'),
    @(20, '11', '72',  '70.15'),
    @(21, '12', '88',  '81.56'),
    @(22, '7',  '49',  '51.85'),
    @(23, '10', '52',  '53.65'),
    @(24, '6',  '26',  '49.95'),
    @(25, '7',  '46',  '89.1'),
    @(26, 'SKIP', '36', '40.42'),
    @(27, '5',  '40',  '51.68'),
    @(28, '7',  '56',  '48.14'),
    @(29, '1',  '4',   '71.43'),
    @(30, '2',  '10',  '46.87'),
    @(31, '1',  '5',   '21.43'),
    @(33, '1',  '3',   '44.05'),
    @(34, '1',  '6',   '54.63'),
    @(35, '1',  '7',   '33.48'),
    @(36, 'SKIP', '15', '15.14'),
    @(37, '1',  '3',   '73.61'),
    @(38, 'CLEAR', 'CLEAR', 'CLEAR'),
    @(39, '1',  '11',  '30.59'),
    @(40, 'CLEAR', 'CLEAR', 'CLEAR'),
    @(41, 'SKIP', '6',  '22.92'),
    @(42, '1',  '7',   '23.47'),
    @(43, 'SKIP', '6',  '12.12')
)

$sheet2Updates = @(
    @(7,  'CLEAR', 'CLEAR', 'CLEAR'),
    @(8,  '2',  '21',  '14.15'),
    @(9,  '4',  '76',  '18.53'),
    @(10, '10', '169', '33.02'),
    @(11, '12', '146', '40.35'),
    @(12, '13', '145', '38.55'),
    @(13, '17', '179', '50.58'),
    @(14, '16', '130', '69.06'),
    @(15, '15', '123', '78.95'),
    @(16, '9',  '65',  '53.03'),
    @(17, '12', '52',  '55.49'),
    @(18, '9',  '68',  '76.41'),
    @(19, '7',  '57',  '43.41'),
    @(20, '7',  '67',  '50.59'),
    @(21, '1',  '7',   '64.58'),
    @(22, '2',  '8',   '25.3'),
    @(24, '1',  '9',   '49.34'),
    @(25, 'SKIP', '7', '33.48'),
    @(26, '2',  '18',  '44.38'),
    @(27, '1',  '5',   '19.5'),
    @(28, '1',  'SKIP', '41.67'),
    @(29, 'SKIP', '6',  '22.92'),
    @(30, '1',  '7',   '23.47'),
    @(31, '1',  '6',   '12.12')
)

function Apply-Updates($ws, $updates) {
    foreach ($u in $updates) {
        $row = $u[0]
        $dVal = $u[1]
        $eVal = $u[2]
        $fVal = $u[3]

        if ($dVal -eq 'CLEAR') {
            $ws.Range("D$row").Value = $null
        } elseif ($dVal -ne 'SKIP') {
            $ws.Range("D$row").Value = [double]$dVal
        }

        if ($eVal -eq 'CLEAR') {
            $ws.Range("E$row").Value = $null
        } elseif ($eVal -ne 'SKIP') {
            $ws.Range("E$row").Value = [double]$eVal
        }

        if ($fVal -eq 'CLEAR') {
            $ws.Range("F$row").Value = $null
        } elseif ($fVal -ne 'SKIP') {
            $ws.Range("F$row").Value = [double]$fVal
        }
    }
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

Apply-Updates $ws1 $sheet1Updates
Apply-Updates $ws2 $sheet2Updates
